$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "25.000.99"
Set-TextValue "E2" "  -3.66%  "
Set-TextValue "D3" "1.637.42"
Set-TextValue "E3" "  -5.62%  "
Set-TextValue "E4" "  -0.18%  "
Set-TextValue "E5" "  -5.75%  "
Set-TextValue "E6" "  -0.11%  "
Set-TextValue "D7" "0.4715"
Set-TextValue "E7" "  -6.00%  "
Set-TextValue "B8" "OKB"
Set-TextValue "C8" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D8" "39.49"
Set-TextValue "E8" "  -3.22%  "
Set-TextValue "B9" "Cardano"
Set-TextValue "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.2548"
Set-TextValue "E9" "  -6.49%  "
Set-TextValue "B10" "Dogecoin"
Set-TextValue "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.06059"
Set-TextValue "E10" "  -1.99%  "
Set-TextValue "B11" "TRON"
Set-TextValue "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D11" "0.07018"
Set-TextValue "E11" "  -3.34%  "
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.641.35"
Set-TextValue "E12" "  -5.47%  "
Set-TextValue "B13" "Solana"
Set-TextValue "C13" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "14.31"
Set-TextValue "E13" "  -5.42%  "
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "4.282"
Set-TextValue "E14" "  -9.77%  "
Set-TextValue "B15" "Polygon"
Set-TextValue "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.5676"
Set-TextValue "E15" "  -13.17%  "
Set-TextValue "B16" "Litecoin"
Set-TextValue "C16" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "73.16"
Set-TextValue "E16" "  -5.71%  "
Set-TextValue "B17" "Dai"
Set-TextValue "C17" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "E17" "  -0.07%  "
Set-TextValue "B18" "BinanceUSD"
Set-TextValue "C18" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D18" "1.000"
Set-TextValue "E18" "  -0.11%  "
Set-TextValue "B19" "WrappedBTC"
Set-TextValue "C19" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D19" "24.989.47"
Set-TextValue "E19" "  -3.76%  "
Set-TextValue "B20" "Avalanche"
Set-TextValue "C20" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D20" "11.18"
Set-TextValue "E20" "  -5.97%  "
Set-TextValue "B21" "ShibaInu"
Set-TextValue "C21" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D21" "0.000006541"
Set-TextValue "E21" "  -4.21%  "
Set-TextValue "B22" "WrappedliquidstakedEther2.0"
Set-TextValue "C22" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D22" "1.851.45"
Set-TextValue "E22" "  -5.87%  "
Set-TextValue "B23" "Uniswap"
Set-TextValue "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D23" "4.250"
Set-TextValue "E23" "  -7.43%  "
Set-TextValue "B24" "Cosmos"
Set-TextValue "C24" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D24" "8.450"
Set-TextValue "E24" "  -3.66%  "
Set-TextValue "B25" "Chainlink"
Set-TextValue "C25" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D25" "5.183"
Set-TextValue "E25" "  -4.24%  "
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "132.30"
Set-TextValue "E26" "  -1.25%  "
Set-TextValue "B27" "EthereumClassic"
Set-TextValue "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D27" "14.79"
Set-TextValue "E27" "  -3.24%  "
Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "1.361"
Set-TextValue "E28" "  -9.72%  "
Set-TextValue "B29" "BitcoinCash"
Set-TextValue "C29" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D29" "103.21"
Set-TextValue "E29" "  -1.99%  "
Set-TextValue "B30" "LidoDAOToken"
Set-TextValue "C30" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D30" "1.624"
Set-TextValue "E30" "  -8.90%  "
Set-TextValue "B31" "InternetComputer(DFINITY)"
Set-TextValue "C31" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "3.866"
Set-TextValue "E31" "  -2.67%  "
Set-TextValue "B32" "Stellar"
Set-TextValue "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D32" "0.07548"
Set-TextValue "E32" "  -7.16%  "
Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "3.502"
Set-TextValue "E33" "  -5.22%  "
Set-TextValue "B34" "Frax"
Set-TextValue "C34" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D34" "0.9997"
Set-TextValue "E34" "  -0.08%  "
Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.04217"
Set-TextValue "E35" "  -10.97%  "
Set-TextValue "B36" "HuobiToken"
Set-TextValue "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.569"
Set-TextValue "E36" "  -3.59%  "
Set-TextValue "B37" "ARBITRUM"
Set-TextValue "C37" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D37" "0.9260"
Set-TextValue "E37" "  -6.96%  "
Set-TextValue "B38" "ImmutableX"
Set-TextValue "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.5870"
Set-TextValue "E38" "  -3.50%  "
Set-TextValue "B39" "MXToken"
Set-TextValue "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D39" "2.564"
Set-TextValue "E39" "  -6.18%  "
Set-TextValue "B40" "TrustWalletToken"
Set-TextValue "C40" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "0.8890"
Set-TextValue "E40" "  +7.16%  "
Set-TextValue "B41" "PaxDollar"
Set-TextValue "C41" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D41" "0.9999"
Set-TextValue "E41" "  -0.09%  "
Set-TextValue "B42" "VeChain"
Set-TextValue "C42" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.01477"
Set-TextValue "E42" "  -8.21%  "
Set-TextValue "B43" "Quant"
Set-TextValue "C43" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D43" "97.57"
Set-TextValue "E43" "  -2.96%  "
Set-TextValue "B44" "RenderToken"
Set-TextValue "C44" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "1.753"
Set-TextValue "E44" "  -9.82%  "
Set-TextValue "B45" "TheSandbox"
Set-TextValue "C45" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D45" "0.3648"
Set-TextValue "E45" "  -6.81%  "
Set-TextValue "B46" "FraxShare"
Set-TextValue "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "4.625"
Set-TextValue "E46" "  -7.62%  "
Set-TextValue "B47" "Algorand"
Set-TextValue "C47" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D47" "0.1089"
Set-TextValue "E47" "  -7.11%  "
Set-TextValue "B48" "Cronos"
Set-TextValue "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.05190"
Set-TextValue "E48" "  -1.75%  "
Set-TextValue "B49" "Aptos"
Set-TextValue "C49" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D49" "6.036"
Set-TextValue "E49" "  -5.02%  "
Set-TextValue "B50" "TrueUSD"
Set-TextValue "C50" "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
Set-TextValue "E50" "  -0.13%  "
Set-TextValue "B51" "USDD"
Set-TextValue "C51" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D51" "1.000"
Set-TextValue "E51" "  -0.07%  "
